# Update the Denmark Summary MSME indicator figures to their more precise
# (rounded-to-hundredths) values. The source cells hold numeric-looking text
# (e.g. "41.8"), so a plain `.Value = "41.79"` assignment would make Excel's
# COM layer auto-coerce the cell to a Number, changing its type. To keep the
# cells as text (matching the original authoring), we briefly mark the cell
# as Text (@) before writing the string, then reset the cell style back to
# "Normal" so no stray formatting/style diff is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $r = $ws.Range($Address)
    $r.NumberFormat = "@"
    $r.Value = $Value
    $r.Style = "Normal"
}

# Enterprises density (per 1000 people) -- Source Type: Statistical Institution
Set-TextValue "B11" "41.79"
Set-TextValue "C11" "8.35"
Set-TextValue "D11" "50.15"

# Source Type: SME Associations (Most Widely Used)
# Enterprises density (per 1000 people)
Set-TextValue "B33" "34.23"
Set-TextValue "C33" "3.96"
Set-TextValue "D33" "38.19"

# Employment (% of total)
Set-TextValue "B34" "21.43"
Set-TextValue "C34" "44.99"
Set-TextValue "D34" "66.43"

# Enterprises (% of total)
Set-TextValue "B36" "89.36"
Set-TextValue "C36" "10.34"
